# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) holds the (wrong) value "4-17-2013-14" for every
# data row; it must become the correct ISO date text "2014-04-17".
#
# Note: the target value looks like a date, so a direct
#   $cell.Value = "2014-04-17"
# would be auto-parsed by Excel into a date serial number instead of being
# kept as literal text. To avoid that (and to avoid introducing any new
# cell formatting/styles), we build the text in a scratch cell via a
# formula that evaluates to the literal string, copy it, and paste only
# the resulting value (PasteSpecial xlPasteValues) into each target cell.
# The scratch row is removed afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "4-17-2013-14"
$newValue = "2014-04-17"

$xlPasteValues = -4163

# Find the last used row in the BF (58th) column.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 58).End(-4162).Row  # xlUp = -4162

# Scratch cell far below the data, used only to generate the literal text
# value without Excel re-interpreting it as a date.
$scratchRow = $ws.Rows.Count
$helper = $ws.Cells.Item($scratchRow, 1)
$helper.Formula = '="' + $newValue + '"'
$helper.Copy()

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 58)
    if ($cell.Value() -eq $oldValue) {
        $cell.PasteSpecial($xlPasteValues)
    }
}

$ws.Rows.Item($scratchRow).Delete()
$excel.CutCopyMode = $false
